# Updated cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape. Column D ("Price") sometimes contains values that read as
# plain numbers (e.g. "1.00", "596.90") -- those must stay as literal text
# (the sheet stores thousands-separated / trailing-zero price strings, not
# numeric values), so every write uses a leading apostrophe to force text
# entry, then resets the cell back to the "Normal" style so no stray
# quote-prefix / number-format style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "67.609.96"
Set-TextValue "E2" "  +1.09%  "
Set-TextValue "D3" "3.490.29"
Set-TextValue "E3" "  -0.20%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "596.90"
Set-TextValue "E5" "  +0.45%  "
Set-TextValue "D6" "179.62"
Set-TextValue "E6" "  +4.00%  "
Set-TextValue "E7" "  +5.70%  "
Set-TextValue "D9" "3.493.01"
Set-TextValue "E9" "  -0.06%  "
Set-TextValue "D10" "0.139"
Set-TextValue "E10" "  +5.36%  "
Set-TextValue "E11" "  -1.45%  "
Set-TextValue "E12" "  +1.36%  "
Set-TextValue "D13" "4.098.96"
Set-TextValue "E13" "  -0.11%  "
Set-TextValue "D14" "32.02"
Set-TextValue "E14" "  +9.37%  "
Set-TextValue "E15" "  +0.07%  "
Set-TextValue "D16" "67.606.22"
Set-TextValue "E16" "  +1.07%  "
Set-TextValue "E17" "  +0.08%  "
Set-TextValue "D18" "3.489.71"
Set-TextValue "E18" "  -1.03%  "
Set-TextValue "D19" "6.32"
Set-TextValue "E19" "  +1.03%  "
Set-TextValue "D20" "14.28"
Set-TextValue "E20" "  -0.24%  "
Set-TextValue "D21" "390.92"
Set-TextValue "E21" "  +0.31%  "
Set-TextValue "D22" "7.97"
Set-TextValue "E22" "  +0.53%  "
Set-TextValue "D23" "73.08"
Set-TextValue "E23" "  -0.33%  "
Set-TextValue "D24" "0.542"
Set-TextValue "E24" "  +1.43%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.48%  "
Set-TextValue "E26" "  +1.32%  "
Set-TextValue "E27" "  +1.44%  "
Set-TextValue "D28" "10.35"
Set-TextValue "E28" "  +2.19%  "
Set-TextValue "E29" "  -2.20%  "
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  +0.24%  "
Set-TextValue "E31" "  +0.39%  "
Set-TextValue "D32" "1.42"
Set-TextValue "E32" "  +0.35%  "
Set-TextValue "E33" "  +0.67%  "
Set-TextValue "D34" "23.59"
Set-TextValue "E34" "  -0.19%  "
Set-TextValue "D35" "7.42"
Set-TextValue "E35" "  +0.68%  "
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  +0.02%  "
Set-TextValue "D37" "1.61"
Set-TextValue "D38" "162.96"
Set-TextValue "E38" "  -0.67%  "
Set-TextValue "E39" "  +0.84%  "
Set-TextValue "D40" "2.83"
Set-TextValue "E40" "  +11.37%  "
Set-TextValue "E41" "  -0.70%  "
Set-TextValue "D42" "6.83"
Set-TextValue "E42" "  -0.30%  "
Set-TextValue "D43" "4.62"
Set-TextValue "E43" "  +0.03%  "
Set-TextValue "D44" "2.831.11"
Set-TextValue "D45" "26.22"
Set-TextValue "E45" "  +1.00%  "
Set-TextValue "D46" "26.78"
Set-TextValue "E46" "  -1.20%  "
Set-TextValue "D47" "0.0727"
Set-TextValue "E47" "  -0.74%  "
Set-TextValue "E48" "  -2.15%  "
Set-TextValue "E49" "  +0.23%  "
Set-TextValue "D50" "335.47"
Set-TextValue "E50" "  -1.34%  "
Set-TextValue "E51" "  -0.98%  "
